$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.153.76"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.838.77"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.01"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6819"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2991"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07453"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.23"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07645"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.838.15"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.026"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6809"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.30"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.156"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -6.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.120.02"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008208"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.071.85"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "230.04"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.34%  "

$ws.Range("E21").Value = "  -2.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.346"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9997"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.96"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.60%  "

$ws.Range("E26").Value = "  -5.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.709"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.61%  "

$ws.Range("E28").Value = "  -1.71%  "

$ws.Range("E29").Value = "  -2.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.253"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.140"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.45%  "

$ws.Range("E32").Value = "  -0.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05340"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7543"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.850"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.132"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.682"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.312.81"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01824"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.718"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9482"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.064"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.69%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "105.16"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9988"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.08141"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +29.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.974.40"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5175"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.77%  "

$ws.Range("E48").Value = "  -3.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.774"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.15"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.391"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.15%  "
